$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J), copying the header style from H1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# row, I value, J value
$rowData = @(
    @(2, 6, 6),
    @(3, 6, 6),
    @(4, 6, 6),
    @(5, 7, 7),
    @(6, 7, 7),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 6, 7),
    @(10, 9, 9),
    @(11, 6, 7),
    @(12, 5, 5),
    @(13, 6, 6),
    @(14, 8, 8),
    @(15, 7, 8),
    @(16, 9, 9),
    @(17, 7, 7),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 8),
    @(21, 9, 10),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 6, 6),
    @(25, 6, 6),
    @(26, 8, 8),
    @(27, 7, 7),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 8, 8),
    @(31, 4, 4),
    @(32, 5, 6),
    @(33, 7, 7),
    @(34, 5, 6),
    @(35, 6, 6),
    @(36, 8, 8),
    @(37, 10, 10),
    @(38, 7, 7),
    @(39, 9, 9),
    @(40, 7, 7),
    @(41, 7, 7),
    @(42, 7, 7),
    @(43, 8, 8),
    @(44, 7, 7),
    @(45, 7, 8),
    @(46, 6, 7),
    @(47, 9, 9),
    @(48, 7, 7),
    @(49, 7, 7),
    @(50, 7, 7),
    @(51, 7, 7),
    @(52, 6, 6),
    @(53, 7, 8),
    @(54, 5, 5),
    @(55, 8, 8),
    @(56, 8, 8),
    @(57, 6, 6),
    @(58, 8, 8),
    @(59, 6, 6),
    @(60, 7, 7),
    @(61, 6, 6),
    @(62, 7, 7),
    @(63, 5, 5),
    @(64, 7, 7),
    @(65, 6, 7),
    @(66, 7, 8),
    @(67, 5, 5),
    @(68, 8, 8),
    @(69, 6, 6),
    @(70, 7, 8),
    @(71, 6, 6),
    @(72, 5, 5),
    @(73, 5, 5),
    @(74, 8, 8),
    @(75, 9, 9),
    @(76, 6, 6),
    @(77, 7, 8),
    @(78, 5, 6),
    @(79, 6, 7),
    @(80, 3, 3)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
